$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full ifoCAST naive-forecaster QoQ error series (rows 2-24, cols B-K).
# The sheet previously held a triangular/staircase matrix (each row had
# fewer trailing columns as the forecast horizon exceeded history available
# at that vintage). With the full ifoCAST series now evaluated, every row
# gets a complete set of Q0-Q9 error values.
$data = @(
    @(-1.282872386776422, 0.1391894435007327, -0.182694870696341, 0.3443354285759472, 0.1070053646347159, 0.2095078597799856, 0.110633009540361, 0.6128907522960124, 0.8708754778845325, -0.8533573809972426),
    @(0.2750606270018259, 0.5534868575902183, 0.2080306981463169, 0.2637694639246414, 0.1446758428348626, 0.6381972808990434, 0.892409660067914, -0.8334509197017121, -0.1806300843769859, 0.1825572214681606),
    @(0.1918994249285815, 0.239828553585737, 0.1203502894047785, 0.6162427231036979, 0.8725296081739512, -0.8520779198862551, -0.1986054693098328, 0.164895676365706, -0.735890762303052, 0.3344822619593895),
    @(0.2023032603298698, 0.6254400633356464, 0.85282204567298, -0.8813162353832781, -0.2307498541254956, 0.1318997099156809, -0.7691308072183836, 0.301173166579823, -0.2759067096457259, -0.4815182333746927),
    @(0.8240207549100882, -0.8834400190618514, -0.2234668893430731, 0.1423162228589251, -0.7576896617840617, 0.3129479666538719, -0.2640230888223336, -0.4695989885331262, 0.5095561729841656, -0.226786726508225),
    @(-0.3453490504388863, 0.05047105392560119, -0.83756240835547, 0.2380186391052039, -0.3369057261414908, -0.5416351540011936, 0.4378696407944773, -0.2983290036150871, -0.3162381895859678, -0.2928084366771486),
    @(-0.7422476366981072, 0.3162863437924887, -0.2663152355556169, -0.4745738317060451, 0.5033030119909454, -0.2336496331081919, -0.2519093219356206, -0.2286430458875994, -0.5938837416006923, 0.06190409654164158),
    @(-0.00619269375287701, -0.3336679001236204, 0.5887475516261532, -0.1739978958108406, -0.2042525916245146, -0.1865645222487989, -0.5543992815454399, 0.1001822482843727, -0.03126473275901087, -0.1410506211186619),
    @(0.3911393015443105, -0.2979577194736586, -0.2956322122566587, -0.2635227406874879, -0.6249641526670209, 0.03245628485901808, -0.09772803603944658, -0.2069514014028122, -1.036602919657539, -0.4560229796881132),
    @(-0.228230158963832, -0.1886450477650728, -0.5467980559473664, 0.1120680736150478, -0.01748087635905388, -0.1264250868990488, -0.9559539928657461, -0.3753202137097185, -0.4586594580766912, -0.6765490042469573),
    @(-0.5591547568613399, 0.09787708368923709, -0.03242713058169561, -0.1416815451927503, -0.9713369016796535, -0.390754202413632, -0.4741138500012969, -0.6920114308521601, 0.8717720887233864, -0.4642089603862817),
    @(0.1398755720077088, -0.04718941243373825, -0.9130114627925084, -0.3492324732466185, -0.4403994753588329, -0.6619245159313032, 0.900173610209962, -0.4365905078372946, 1.003661313589152, -0.09695126965353879),
    @(-1.007520009675017, -0.3857468722612991, -0.4502017751758141, -0.6594324126472006, 0.9083245926860781, -0.4258347558073555, 1.015616079910623, -0.0844445585108285, -0.2945431432588089, 0.3040816658791113),
    @(-0.3077970529158559, -0.5788562590806742, 0.961338929688949, -0.3851068719675381, 1.050866189894518, -0.05163699300024938, -0.2628248666403167, 0.3353140868060093, 0.009465662865194846, 0.7265657060874963),
    @(1.139271322614915, -0.2982492763927073, 1.093791915679716, -0.02972831832925377, -0.250970648745674, 0.3423607205295725, 0.01421452142035373, 0.7302167878130713, 2.425529404116681, 9.269510203911928),
    @(-0.2876812562600493, 1.102566704251488, -0.02192591428551849, -0.2436605879111656, 0.3494215681832005, 0.02114969792139543, 0.7370887910462527, 2.432369740824074, 9.276334707790259, -8.274453695494744),
    @(0.9955923004358076, -0.0871855372237873, -0.2892511864441449, 0.3130908916669468, -0.0108203968125144, 0.7071725634927386, 2.403421102624959, 9.247842001542866, -8.302731520459806, -0.5932976437114486),
    @(-0.0533123394792912, -0.2864980261873478, 0.3015019608161922, -0.02894181086942849, 0.6860724501299835, 2.380962765938485, 9.224764288378545, -8.32609170460871, -0.6167866614867294, 0.9660844828649828),
    @(-0.5437734556935978, 0.1702226893339956, -0.1022974921431762, 0.63883117546281, 2.345515881559747, 9.194643401784708, -8.35380748862633, -0.6434163542466904, 0.9399452457891766, -2.026139092643445),
    @(0.1269585753652727, -0.1306388881700511, 0.6141551499506877, 2.321229248582915, 9.170020105529264, -8.378790661567026, -0.668646313777596, 0.9145696739774734, -2.051594272676513, 0.2124867959412257),
    @(-0.02044773294388247, 0.6702777838950345, 2.355311064516826, 9.195119489401716, -8.357352010851283, -0.6486996043872704, 0.933908325411367, -2.032503447757091, 0.2314766118486793, -0.2196984654044191),
    @(0.5134276986299353, 2.260143333751186, 9.140692358735786, -8.39390599826305, -0.6763647928704721, 0.9105011857394278, -2.053840367633118, 0.2111440141036455, -0.2395424444363658, -0.379677132609091),
    @(2.211279720138829, 9.112923508459115, -8.41112584330054, -0.688459908824067, 0.9009170582553292, -2.06219729845007, 0.2033872967055401, -0.2470056675122654, -0.386996833524363, 0.1394248687261353)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $excelCol = $j + 2
        $ws.Cells.Item($excelRow, $excelCol).Value2 = $rowValues[$j]
    }
}
